$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 becomes a literal number instead of the shared string "test"
$ws.Range("A2").Value = 123456

# C2 value changes from 100 to 123456789
$ws.Range("C2").Value = 123456789

# Column C gets an explicit width (best-fit) matching the new wider content
$ws.Columns.Item(3).AutoFit() | Out-Null
